# ---------------------------------------------------------------------------
# C1--C2-and-C3-PowerPoint.pptx edit
#
# 1) Slide 16: the table on the slide switches from the custom "Table_0"
#    style to the built-in "No Style, Table Grid" style
#    ({28912C55-EE52-4C6D-A18E-44F83DEEE6AA}).
#
# 2) The deck's theme colour scheme (used by the slide master / all slides)
#    is switched from the "Integral" palette back to the default "Office"
#    palette (i.e. the Design/Theme was changed back to the stock Office
#    theme).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 --------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{28912C55-EE52-4C6D-A18E-44F83DEEE6AA}")
    }
}

# --- 2) Theme colours: Integral -> Office -----------------------------------
$master = $p.Designs.Item(1).SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
